# Consolidate the multiple <a:r> text runs on slide 1 into single runs,
# the same way PowerPoint does when it re-writes a paragraph's runs that
# all share identical run properties (<a:rPr/>).
#
# Plain `TextRange.Text = "..."` only rewrites the *first* run of a
# paragraph and leaves the remaining runs untouched, so instead we target
# the exact character ranges (via `Characters(Start, Length)`) that make
# up each paragraph "segment" (i.e. the text between two manual line
# breaks) and replace it in one shot; this merges all the runs inside
# that segment into a single run while leaving the <a:br/> elements
# between segments alone.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1 ("Title 1"): "Testing" + " " + "custom" + " " + "properties"
#     -> single run "Testing custom properties"
$trTitle = $s.Shapes.Item(1).TextFrame.TextRange
$trTitle.Characters(1, $trTitle.Length).Text = "Testing custom properties"

# --- Shape 2 ("Subtitle 2"):
#     "This" " " "is" " " "a" " " "subtitle" <br/> <br/> "A." " " "M."
#     -> "This is a subtitle" <br/> <br/> "A. M."
$trSub = $s.Shapes.Item(2).TextFrame.TextRange
$trSub.Characters(1, 18).Text = "This is a subtitle"
$trSub.Characters(21, 5).Text = "A. M."
